# Update latest output (run 28)
# Applies the refreshed optimisation_result numbers to the "Schedule" and
# "Detailed" worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Schedule": refreshed cost / unit-cost figures (columns E & F)
# ---------------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("E2").Value = 1178.2416555
$schedule.Range("F2").Value = 25.97534513888889

$schedule.Range("E3").Value = 354.6176205
$schedule.Range("F3").Value = 23.45354632936508

# ---------------------------------------------------------------------
# Sheet "Detailed": refreshed Price (col B) and, for a few early-morning
# rows, the Type flips from "forecast" to "historical" (col C) as the
# actual data caught up to those timestamps.
# ---------------------------------------------------------------------
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B5").Value = 78
$detailed.Range("B6").Value = 78

$detailed.Range("B7").Value = 87.05036
$detailed.Range("C7").Value = "historical"

$detailed.Range("C8").Value = "historical"

$detailed.Range("B9").Value = 85.44701000000001
$detailed.Range("C9").Value = "historical"

$detailed.Range("B11").Value = 80.02
$detailed.Range("B12").Value = 80.02
$detailed.Range("B13").Value = 93.26281

$detailed.Range("B16").Value = 56.98
$detailed.Range("B18").Value = 56.97999
$detailed.Range("B19").Value = 56.90274
$detailed.Range("B20").Value = 47.13666
$detailed.Range("B22").Value = 46.68963
$detailed.Range("B23").Value = 41.21992
$detailed.Range("B24").Value = 36.07
$detailed.Range("B27").Value = 36.07

$detailed.Range("B32").Value = 27.27348
$detailed.Range("B33").Value = 18.62722
$detailed.Range("B34").Value = 17.79394
$detailed.Range("B35").Value = 7.7265
$detailed.Range("B36").Value = -6

$detailed.Range("B38").Value = -3.13002
$detailed.Range("B39").Value = -3.05272
$detailed.Range("B40").Value = 0.01129

$detailed.Range("B44").Value = 0.08645
$detailed.Range("B45").Value = 73.73759
$detailed.Range("B46").Value = 57.55625
$detailed.Range("B47").Value = 57.44178
$detailed.Range("B48").Value = 58.20187
